# Pax and Crew Manifest Update
# - Replace the "Age" column header with "Date Of Birth(MM/DD/YYYY)" (column L)
# - Row 1 height grows to fit the new, longer wrapped header text
# - Selection moves to the newly edited header cell (L1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header text in column L ("Age" -> "Date Of Birth(MM/DD/YYYY)")
$ws.Range("L1").Value = "Date Of Birth(MM/DD/YYYY)"

# Row 1 needs to be taller to accommodate the new wrapped text
$ws.Rows.Item(1).RowHeight = 60

# Reflect the new selection/active cell state
$ws.Range("L1").Select()

$wb.Save()
